# Update NATMI Wnt1-Ror2 LR-pair output after recomputing TPM without the
# "MuSCs" cluster: drop the rows that reference MuSCs (old rows 6 & 7, plus
# the "MuSCs" target-cluster row among ECs/FAPs), and refresh the derived
# statistics for the remaining ECs/FAPs combinations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two bottom rows that belonged to the FAPs->FAPs / FAPs->MuSCs
# combinations that are being replaced/removed; remaining rows shift up.
$ws.Range("A6:T7").EntireRow.Delete()

# --- Row 2: ECs -> ECs (unchanged pairing, refreshed metrics) ---
$ws.Range("G2").Value = 0.027123
$ws.Range("H2").Value = 0.081369
$ws.Range("I2").Value = 0.0960827240265261
$ws.Range("J2").Value = 0.09608272402652611
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1074926666666667
$ws.Range("N2").Value = 0.322478
$ws.Range("O2").Value = 0.01930181557781338
$ws.Range("P2").Value = 0.01930181557781338
$ws.Range("Q2").Value = 0.002915523598
$ws.Range("R2").Value = 0.026239712382
$ws.Range("S2").Value = 0.001854571019373946
$ws.Range("T2").Value = 0.001854571019373946

# --- Row 3: ECs -> FAPs (unchanged pairing, refreshed metrics) ---
$ws.Range("G3").Value = 0.027123
$ws.Range("H3").Value = 0.081369
$ws.Range("I3").Value = 0.0960827240265261
$ws.Range("J3").Value = 0.09608272402652611
$ws.Range("M3").Value = 5.461551666666666
$ws.Range("O3").Value = 0.9806981844221867
$ws.Range("P3").Value = 0.9806981844221866
$ws.Range("Q3").Value = 0.148133665855
$ws.Range("R3").Value = 1.333202992695
$ws.Range("S3").Value = 0.09422815300715216
$ws.Range("T3").Value = 0.09422815300715216

# --- Row 4: now FAPs -> ECs (previously FAPs -> MuSCs; MuSCs dropped) ---
$ws.Range("A4").Value = "FAPs"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.255165
$ws.Range("H4").Value = 0.7654949999999999
$ws.Range("I4").Value = 0.9039172759734738
$ws.Range("J4").Value = 0.9039172759734738
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.1074926666666667
$ws.Range("N4").Value = 0.322478
$ws.Range("O4").Value = 0.01930181557781338
$ws.Range("P4").Value = 0.01930181557781338
$ws.Range("Q4").Value = 0.02742836629
$ws.Range("R4").Value = 0.24685529661
$ws.Range("S4").Value = 0.01744724455843944
$ws.Range("T4").Value = 0.01744724455843943

# --- Row 5: now FAPs -> FAPs (previously FAPs -> ECs) ---
$ws.Range("D5").Value = "FAPs"
$ws.Range("I5").Value = 0.9039172759734738
$ws.Range("J5").Value = 0.9039172759734738
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 5.461551666666666
$ws.Range("N5").Value = 16.384655
$ws.Range("O5").Value = 0.9806981844221867
$ws.Range("P5").Value = 0.9806981844221866
$ws.Range("Q5").Value = 1.393596831025
$ws.Range("R5").Value = 12.542371479225
$ws.Range("S5").Value = 0.8864700314150344
$ws.Range("T5").Value = 0.8864700314150343
